$d = $word.ActiveDocument

$replacements = @(
    @{old = "80×29="; new = "18×92="},
    @{old = "25×22="; new = "70×66="},
    @{old = "89×97="; new = "70×67="},
    @{old = "41×43="; new = "15×13="},
    @{old = "88×94="; new = "34×65="},
    @{old = "12×49="; new = "50×75="},
    @{old = "45×66="; new = "62×57="},
    @{old = "43×57="; new = "45×71="},
    @{old = "39×52="; new = "55×64="},
    @{old = "98×71="; new = "52×76="},
    @{old = "40×35="; new = "70×80="},
    @{old = "69×68="; new = "82×59="},
    @{old = "65×98="; new = "90×94="},
    @{old = "12×42="; new = "67×37="},
    @{old = "11×30="; new = "48×40="},
    @{old = "40×50="; new = "24×44="},
    @{old = "20×47="; new = "34×67="},
    @{old = "50×42="; new = "92×14="},
    @{old = "65×84="; new = "75×13="},
    @{old = "77×14="; new = "44×19="},
    @{old = "24×47="; new = "13×64="},
    @{old = "75×63="; new = "59×46="},
    @{old = "56×49="; new = "73×48="},
    @{old = "56×17="; new = "75×39="},
    @{old = "36×30="; new = "90×20="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
